# The commit updates the "Förändrad" (Changed) date column (C) for every
# data row (2-205) from serial date 45182 (2023-09-13) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 205; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
